# Regenerate save_data to use K (strikeouts) instead of Strike# and
# recompute the K column (G) values for each start.
#
# New K values per row (A2:A35 correspond to game index 0..33),
# worksheet rows 2..35. Row 34 is unchanged (already 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,0,2,0,3,2,3,2,5,0,0,1,4,2,1,2,8,6,1,6,2,0,6,4,5,2,4,4,3,1,2,3,0,3)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
